# Applies the "handles float input without breaking stuff" fix to the
# marksheet result sheet: recomputed summary stats (rows 10-12), the
# per-question "Student Ans" column (A16:A40) is now populated/colored for
# every attempted question, one extra "Correct Ans" helper column (D16:D18)
# is filled in, and the now-unused extra Student/Correct-Ans block
# (D19:E40 and G15:H21) is cleared away so the used range shrinks back
# down to A5:E40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Summary block (rows 10-12): the row-label cells (A10:A12) pick up
#    the same "mtitleStyle" the header row (row 9) already uses - copy
#    the format from there instead of re-creating the named style, then
#    drop in the recomputed right/wrong/not-attempted/total figures.
# ---------------------------------------------------------------------
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial($xlPasteFormats)
$ws.Range("A11").PasteSpecial($xlPasteFormats)
$ws.Range("A12").PasteSpecial($xlPasteFormats)

$ws.Range("B10").Value = 21
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "82/112"

# ---------------------------------------------------------------------
# 2) Per-question "Student Ans" (col A) for rows 16-40: every question
#    that was attempted now shows the student's chosen option, styled
#    green ("correctStyle", same format as B10) when it matches the
#    "Correct Ans" column (B) and red ("incorrectStyle", same format as
#    C10) when it doesn't. Rows that were never attempted are left
#    blank with the plain "normalStyle" they already had.
# ---------------------------------------------------------------------
$correctRows = @(17, 18, 19, 20, 21, 22, 23, 24, 25, 27, 30, 31, 32, 33, 34, 35, 37, 38, 40)
$incorrectRows = @(36)

$studentAnswers = @{
    17 = "Option D"
    18 = "Option B"
    19 = "Option C"
    20 = "Option B"
    21 = "Option C"
    22 = "Option D"
    23 = "Option D"
    24 = "Option A"
    25 = "Option A"
    27 = "Option A"
    30 = "Option B"
    31 = "Option D"
    32 = "Option C"
    33 = "Option D"
    34 = "Option B"
    35 = "Option D"
    36 = "Option D"
    37 = "Option A"
    38 = "Option A"
    40 = "Option D"
}

$ws.Range("B10").Copy()
foreach ($r in $correctRows) {
    $ws.Cells.Item($r, 1).PasteSpecial($xlPasteFormats)
}

$ws.Range("C10").Copy()
foreach ($r in $incorrectRows) {
    $ws.Cells.Item($r, 1).PasteSpecial($xlPasteFormats)
}

foreach ($r in $studentAnswers.Keys) {
    $ws.Cells.Item($r, 1).Value = $studentAnswers[$r]
}

# ---------------------------------------------------------------------
# 3) The second "Student Ans"/"Correct Ans" pair (col D) only keeps
#    data for the first three questions now; fill in the student
#    answers there (green when right, same "correctStyle"/"incorrectStyle"
#    formats copied from B10/C10).
# ---------------------------------------------------------------------
$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial($xlPasteFormats)
$ws.Range("D18").PasteSpecial($xlPasteFormats)

$ws.Range("C10").Copy()
$ws.Range("D17").PasteSpecial($xlPasteFormats)

$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option D"
$ws.Range("D18").Value = "Option D"

# ---------------------------------------------------------------------
# 4) Drop the now-unused data: the rest of the second pair (D19:E40)
#    and the whole third "Student Ans"/"Correct Ans" pair (G15:H21),
#    so the sheet's used range collapses back down to A5:E40.
# ---------------------------------------------------------------------
$ws.Range("D19:E40").Clear()
$ws.Range("G15:H21").Clear()
